# Commit: "Tweak jpg for better aspect ratio on Mixer"
#
# On the Mixer Stream Schedule slide (slide 5), two pictures are
# repositioned/resized:
#   - "Picture 20" (the Digital Skills All-Star star/banner graphic):
#       moved down (Top only changes; Left/Width/Height stay the same).
#   - "Picture 2" (the cropped photo, desc "...shirt..."): repositioned
#       and enlarged (Left/Top/Width/Height all change), giving it a
#       better aspect ratio.
#
# NOTE: Left/Top/Width/Height are expressed in points in the PowerPoint
# object model (1 pt = 12700 EMU), and are stored internally with
# single-precision (float32) granularity, so values below are chosen to
# round-trip to the exact target EMU values from the canonical OOXML
# after PowerPoint's internal float32 storage.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# "Picture 20" -- star/banner graphic: only the vertical position moves
# (x=7881105 unchanged; y: 2475425 -> 2680146 EMU).
$picStar = $s.Shapes.Item(5)
$picStar.Top = 211.03512573242188

# "Picture 2" -- shirt/photo cutout: reposition and resize
# (x: 2198323 -> 2118957, y: 3210745 -> 3952274,
#  cx: 810762 -> 919944, cy: 1534346 -> 1740970 EMU).
$picShirt = $s.Shapes.Item(7)
$picShirt.Left = 166.84701538085938
$picShirt.Top = 311.20269775390625
$picShirt.Width = 72.43653869628906
$picShirt.Height = 137.08425903320312
